$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Set columns E through AC (years 2010-2034) to 0 for rows 15, 40 and 65
# (id_heating_technology = 213), leaving AD onward (2035+) untouched.
$rows = @(15, 40, 65)
foreach ($r in $rows) {
    $ws.Range("E${r}:AC${r}").Value = 0
}

# Update the selected cell/range shown when the sheet is reopened.
$ws.Range("AG87").Select()

$wb.Save()
